# Updates cryptos list values (price/volume/coin order) per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.348.17'
$ws.Range("E2").Value = '  +0.16%  '

# Row 3
$ws.Range("D3").Value = '2.648.84'
$ws.Range("E3").Value = '  +0.30%  '

# Row 4
$ws.Range("E4").Value = '  +0.22%  '

# Row 5
$ws.Range("D5").Value = "'595.92"
$ws.Range("E5").Value = '  -0.38%  '

# Row 6
$ws.Range("D6").Value = "'158.50"
$ws.Range("E6").Value = '  +2.52%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").Value = "'0.539"
$ws.Range("E8").Value = '  -0.94%  '

# Row 9
$ws.Range("D9").Value = '2.649.22'
$ws.Range("E9").Value = '  +0.36%  '

# Row 10
$ws.Range("D10").Value = "'0.139"
$ws.Range("E10").Value = '  -4.17%  '

# Row 11
$ws.Range("E11").Value = '  -0.85%  '

# Row 12
$ws.Range("D12").Value = "'5.26"
$ws.Range("E12").Value = '  +0.28%  '

# Row 13
$ws.Range("E13").Value = '  +0.41%  '

# Row 14
$ws.Range("D14").Value = "'27.96"
$ws.Range("E14").Value = '  +0.18%  '

# Row 15
$ws.Range("D15").Value = '3.133.06'
$ws.Range("E15").Value = '  +0.35%  '

# Row 16
$ws.Range("D16").Value = "'0.0000186"
$ws.Range("E16").Value = '  -3.86%  '

# Row 17
$ws.Range("D17").Value = '68.273.45'
$ws.Range("E17").Value = '  +0.04%  '

# Row 18
$ws.Range("D18").Value = '2.672.50'
$ws.Range("E18").Value = '  +1.25%  '

# Row 19
$ws.Range("D19").Value = "'11.63"
$ws.Range("E19").Value = '  +2.44%  '

# Row 20
$ws.Range("D20").Value = "'363.34"
$ws.Range("E20").Value = '  -0.07%  '

# Row 21
$ws.Range("D21").Value = "'7.49"
$ws.Range("E21").Value = '  +1.13%  '

# Row 22
$ws.Range("D22").Value = "'4.43"
$ws.Range("E22").Value = '  +1.12%  '

# Row 23
$ws.Range("B23").Value = 'NEARProtocol'
$ws.Range("C23").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D23").Value = "'4.78"
$ws.Range("E23").Value = '  -0.99%  '

# Row 24
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").Value = "'2.09"
$ws.Range("E24").Value = '  +1.38%  '

# Row 25
$ws.Range("D25").Value = "'75.00"
$ws.Range("E25").Value = '  +0.34%  '

# Row 26
$ws.Range("E26").Value = '  +0.05%  '

# Row 27
$ws.Range("D27").Value = "'9.93"
$ws.Range("E27").Value = '  +1.47%  '

# Row 28
$ws.Range("D28").Value = '2.787.42'
$ws.Range("E28").Value = '  +0.51%  '

# Row 29
$ws.Range("D29").Value = "'0.0000102"
$ws.Range("E29").Value = '  -3.68%  '

# Row 30
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = '  +0.20%  '

# Row 31
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = "'577.00"
$ws.Range("E31").Value = '  +2.48%  '

# Row 32
$ws.Range("D32").Value = "'8.13"
$ws.Range("E32").Value = '  +1.46%  '

# Row 33
$ws.Range("D33").Value = "'1.39"
$ws.Range("E33").Value = '  -0.46%  '

# Row 34
$ws.Range("D34").Value = "'1.88"
$ws.Range("E34").Value = '  +0.87%  '

# Row 35
$ws.Range("E35").Value = '  +3.37%  '

# Row 36
$ws.Range("E36").Value = '  +0.03%  '

# Row 37
$ws.Range("E37").Value = '  +0.13%  '

# Row 38
$ws.Range("D38").Value = "'160.55"
$ws.Range("E38").Value = '  -0.07%  '

# Row 39
$ws.Range("D39").Value = "'19.60"
$ws.Range("E39").Value = '  +1.52%  '

# Row 40
$ws.Range("E40").Value = '  -0.81%  '

# Row 41
$ws.Range("E41").Value = '  +0.04%  '

# Row 42
$ws.Range("D42").Value = "'5.31"
$ws.Range("E42").Value = '  -0.19%  '

# Row 43
$ws.Range("D43").Value = "'2.62"
$ws.Range("E43").Value = '  -0.91%  '

# Row 44
$ws.Range("D44").Value = '0.0₆0317'
$ws.Range("E44").Value = '  -6.33%  '

# Row 45
$ws.Range("E45").Value = '  +0.04%  '

# Row 46
$ws.Range("D46").Value = "'158.31"
$ws.Range("E46").Value = '  +0.18%  '

# Row 47
$ws.Range("D47").Value = "'3.82"
$ws.Range("E47").Value = '  +2.05%  '

# Row 48
$ws.Range("D48").Value = "'1.73"
$ws.Range("E48").Value = '  +1.93%  '

# Row 49
$ws.Range("D49").Value = "'21.71"
$ws.Range("E49").Value = '  -0.86%  '

# Row 50
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").Value = "'0.585"
$ws.Range("E50").Value = '  +4.48%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.0780"
$ws.Range("E51").Value = '  -0.73%  '
